$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text
        if ($t -eq "A") {
            $sh.TextFrame.TextRange.Text = "a"
        } elseif ($t -eq "B") {
            $sh.TextFrame.TextRange.Text = "b"
        } elseif ($t -eq "E") {
            $sh.TextFrame.TextRange.Text = "e"
        } elseif ($t -eq "D") {
            $sh.TextFrame.TextRange.Text = "d"
        } elseif ($t -eq "C") {
            $sh.TextFrame.TextRange.Text = "c"
        }
    }
}
